$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.2
$ws.Range("H2").Value = 2.2
$ws.Range("K2").Value = 3.55
$ws.Range("L2").Value = 1.44
$ws.Range("V2").Value = 1.8
$ws.Range("W2").Value = 1.32
$ws.Range("Y2").Value = 9.6
$ws.Range("AA2").Value = 32
$ws.Range("AB2").Value = 15
$ws.Range("AE2").Value = 29
$ws.Range("AG2").Value = 19
$ws.Range("AO2").Value = 24
$ws.Range("F3").Value = 1.86
$ws.Range("I3").Value = 5.1
$ws.Range("K3").Value = 3.9
$ws.Range("N3").Value = 3.55
$ws.Range("O3").Value = 1.34
$ws.Range("P3").Value = 1.87
$ws.Range("Q3").Value = 1.96
$ws.Range("S3").Value = 3.55
$ws.Range("T3").Value = 1.81
$ws.Range("U3").Value = 2
$ws.Range("X3").Value = 17.5
$ws.Range("Y3").Value = 19
$ws.Range("Z3").Value = 42
$ws.Range("AC3").Value = 9.800000000000001
$ws.Range("AH3").Value = 21
$ws.Range("AN3").Value = 16
$ws.Range("F4").Value = 2.64
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 2.74
$ws.Range("I4").Value = 3.15
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 3.55
$ws.Range("M4").Value = 1.05
$ws.Range("V4").Value = 1.46
$ws.Range("W4").Value = 1.5
$ws.Range("F5").Value = 1.49
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 8.800000000000001
$ws.Range("J5").Value = 4.1
$ws.Range("K5").Value = 5
$ws.Range("L5").Value = 1.36
$ws.Range("N5").Value = 3.65
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 1.93
$ws.Range("Q5").Value = 1.86
$ws.Range("R5").Value = 1.36
$ws.Range("U5").Value = 1.83
$ws.Range("W5").Value = 2.72
$ws.Range("AB5").Value = 1000
$ws.Range("AF5").Value = 9.4
$ws.Range("F6").Value = 3.15
$ws.Range("G6").Value = 3.45
$ws.Range("H6").Value = 2.36
$ws.Range("I6").Value = 2.5
$ws.Range("K6").Value = 3.7
$ws.Range("L6").Value = 1.39
$ws.Range("O6").Value = 1.32
$ws.Range("P6").Value = 1.94
$ws.Range("R6").Value = 1.35
$ws.Range("S6").Value = 3.35
$ws.Range("T6").Value = 1.73
$ws.Range("U6").Value = 2.26
$ws.Range("AB6").Value = 13
$ws.Range("AC6").Value = 8.199999999999999
$ws.Range("AF6").Value = 25
$ws.Range("AH6").Value = 17
$ws.Range("F7").Value = 1.39
$ws.Range("H7").Value = 8.800000000000001
$ws.Range("I7").Value = 10.5
$ws.Range("J7").Value = 4.8
$ws.Range("L7").Value = 1.22
$ws.Range("O7").Value = 1.21
$ws.Range("Q7").Value = 1.64
$ws.Range("T7").Value = 1.87
$ws.Range("W7").Value = 3.2
$ws.Range("AK7").Value = 18
$ws.Range("AM7").Value = 150
$ws.Range("AN7").Value = 6
$ws.Range("F8").Value = 1.65
$ws.Range("G8").Value = 1.74
$ws.Range("H8").Value = 6.8
$ws.Range("I8").Value = 8.800000000000001
$ws.Range("J8").Value = 3.4
$ws.Range("L8").Value = 1.54
$ws.Range("N8").Value = 2.56
$ws.Range("O8").Value = 1.57
$ws.Range("P8").Value = 1.51
$ws.Range("R8").Value = 1.18
$ws.Range("S8").Value = 4.9
$ws.Range("U8").Value = 1.58
$ws.Range("W8").Value = 2.34
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 2.32
$ws.Range("H9").Value = 3.4
$ws.Range("N9").Value = 3.7
$ws.Range("P9").Value = 1.91
$ws.Range("Q9").Value = 1.78
$ws.Range("S9").Value = 3
$ws.Range("T9").Value = 1.68
$ws.Range("V9").Value = 1.32
$ws.Range("W9").Value = 1.76
